$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two trailing duplicate rows (old rows 59 and 60)
$ws.Range("A60").EntireRow.Delete()
$ws.Range("A59").EntireRow.Delete()

# Update cell values that changed between before and after
$ws.Range("C4").Value = "str14"

$ws.Range("B29").Value = "CRP_nr"
$ws.Range("C29").Value = "double"
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = "2010USD avg per-CRPacre contract-based FY rent payments (not actuals) (USDA)"

$ws.Range("B30").Value = "crop_nr"
$ws.Range("C30").Value = "double"
$ws.Range("E30").Value = "2010USD annualized net return/acre net income deriving from crop production[L&M]"

$ws.Range("B31").Value = "forest_nr"
$ws.Range("E31").Value = "2010USD annualized net return/acre of bare forestland [L&M]"

$ws.Range("B32").Value = "urban_nr"
$ws.Range("E32").Value = "2010USD annualized net return/acre derived from price of recently dev. land[L&M]"

$ws.Range("B33").Value = "range_nr"
$ws.Range("C33").Value = "float"
$ws.Range("E33").Value = "= pasture_nr"

$ws.Range("B34").Value = "pasture_nr_level"
$ws.Range("C34").Value = "str17"
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = "Pasture rents (NASS) data level"

$ws.Range("B35").Value = "CRPacresk"
$ws.Range("C35").Value = "float"
$ws.Range("E35").Value = "Thousand acres in CRP (USDA County Stats)"

$ws.Range("B36").Value = "CRPrent"
$ws.Range("C36").Value = "double"
$ws.Range("D36").Value = 1
$ws.Range("E36").Value = "2010USD CRP Contract-based FY rent payments (not actuals) (USDA)"

$ws.Range("B37").Value = "CRPland_acresk"
$ws.Range("C37").Value = "float"
$ws.Range("D37").Value = 1

$ws.Range("B38").Value = "Cropland_acresk"
$ws.Range("C38").Value = "float"
$ws.Range("D38").Value = 1

$ws.Range("B39").Value = "Federalland_acresk"

$ws.Range("B40").Value = "Forestland_acresk"

$ws.Range("B41").Value = "Pastureland_acresk"

$ws.Range("B42").Value = "Rangeland_acresk"

$ws.Range("B43").Value = "Ruralland_acresk"

$ws.Range("B44").Value = "Urbanland_acresk"

$ws.Range("B45").Value = "Waterland_acresk"

$ws.Range("B46").Value = "lccNA_acresk"

$ws.Range("B47").Value = "lccL1_acresk"

$ws.Range("B48").Value = "lccL2_acresk"

$ws.Range("B49").Value = "lccL3_acresk"

$ws.Range("B50").Value = "lccL4_acresk"

$ws.Range("B51").Value = "lccL5_acresk"

$ws.Range("B52").Value = "lccL6_acresk"

$ws.Range("B53").Value = "lccL7_acresk"

$ws.Range("B54").Value = "lccL8_acresk"

$ws.Range("B55").Value = "lccL12_acresk"

$ws.Range("B56").Value = "lccL34_acresk"

$ws.Range("B57").Value = "lccL56_acresk"

$ws.Range("B58").Value = "lccL78_acresk"

